# Update "想去人数" (want-to-go count, column F) values for several
# events across the "展览" and "全部类型" worksheets to match the
# newly scraped output (commit: "Update gh-pages to output generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# Sheet "展览" (first sheet)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F14").Value = 47
$wsExhibit.Range("F18").Value = 2463
$wsExhibit.Range("F26").Value = 487
$wsExhibit.Range("F28").Value = 2485
$wsExhibit.Range("F32").Value = 182
$wsExhibit.Range("F40").Value = 2232

# Sheet "全部类型" (fourth sheet, aggregated view with the same events)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F15").Value = 47
$wsAll.Range("F19").Value = 2463
$wsAll.Range("F28").Value = 487
$wsAll.Range("F30").Value = 2485
$wsAll.Range("F34").Value = 182
$wsAll.Range("F45").Value = 2232

$wb.Save()
